# Add "Best exam" and "Worst exam" rows (14 and 15) to each worksheet,
# with the per-sheet NUM_OF_EXAMINED_STM values in column C.

$wb = $excel.ActiveWorkbook

$examValues = @{
    "Tarantula"       = @(2.618726302936831, 16.1461345671872)
    "Ochiai"          = @(1.423430370798795, 18.48210269262902)
    "Op2"             = @(3.85900649058544, 24.84737484737486)
    "Barinel"         = @(2.618726302936831, 16.1461345671872)
    "Dstar"           = @(1.397725081935611, 21.40607930081614)
    "Russell_rao"     = @(7.075380759591283, 30.77244393033867)
    "Simple_matching" = @(13.61095045305574, 45.6622325043378)
    "Rogers_tanimoto" = @(13.61095045305574, 45.6622325043378)
    "Ample"           = @(1.452348820769875, 34.57361352098194)
    "Jaccard"         = @(2.13996529786004, 18.8162714478504)
    "Cohen"           = @(2.229933808881182, 17.29323308270676)
    "Scott"           = @(7.393483709273198, 35.63395668658827)
    "Rogot1"          = @(7.393483709273198, 35.63395668658827)
    "Geometric_mean"  = @(1.471627787417265, 16.1364950838635)
    "M2"              = @(1.645138487243752, 24.39431913116124)
    "Wong1"           = @(7.075380759591283, 30.77244393033867)
    "Sokal"           = @(13.61095045305574, 45.6622325043378)
    "Sorensen_dice"   = @(2.13996529786004, 18.8162714478504)
    "Dice"            = @(2.13996529786004, 18.8162714478504)
    "Humman"          = @(13.61095045305574, 45.6622325043378)
    "Wong2"           = @(13.61095045305574, 45.6622325043378)
    "Euclid"          = @(13.61095045305574, 45.6622325043378)
    "Zoltar"          = @(1.908617698091382, 19.17614549193497)
    "Rogot2"          = @(1.429856693014591, 16.70201143885354)
    "Hamming"         = @(13.61095045305574, 45.6622325043378)
    "Fleiss"          = @(7.734078786710384, 37.94100636205899)
    "Anderberg"       = @(2.13996529786004, 18.8162714478504)
    "Goodman"         = @(2.13996529786004, 18.8162714478504)
    "Harmonic_mean"   = @(1.429856693014591, 16.70201143885354)
    "Kulczynski2"     = @(1.523038365143631, 20.13045434098066)
}

foreach ($ws in $wb.Worksheets) {
    $vals = $examValues[$ws.Name]
    if ($vals -eq $null) { continue }

    $ws.Range("A14").Value = "Best exam"
    $ws.Range("C14").Value = $vals[0]

    $ws.Range("A15").Value = "Worst exam"
    $ws.Range("C15").Value = $vals[1]
}
